# Add an intro schema
#
# Paragraph 5.2 (the "list of figures" entry) originally read:
#   "      5.2 Кроссплатформенное веб-приложение ведения бюджета "il
#   budgetto"». Схема структурная. "
# It is being rewritten to:
#   "      5.2 Структура приложения. Схема структурная."
# keeping the single red-colored space that sits between the two
# sentences, and dropping the trailing bold space that used to close
# the paragraph.

$d = $word.ActiveDocument

# Locate the target paragraph robustly (rather than relying on a fixed
# character offset) by scanning for the one that still has the old
# wording.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -like "*budgetto*") -and ($t -like "*структурная*")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pr = $target.Range
    $q = [char]34
    $guillemet = [char]187

    # Replace the whole old sentence (through the trailing bold space)
    # with the new plain sentence. Starting the match right at
    # "Кроссплатформенное" and ending it past the final space removes
    # the obsolete spell-check markers together with the old runs, and
    # also removes the stray trailing bold space run entirely.
    $find = "Кроссплатформенное веб-приложение ведения бюджета " + $q + "il budgetto" + $q + $guillemet + ". Схема структурная. "
    $replace = "Структура приложения. Схема структурная."
    [void]$pr.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)

    # Touch across the "5.2 " / "Структура" boundary (replacing it with
    # itself) so any now-orphaned proofing-error marker there gets
    # dropped during re-serialization.
    $rng = $d.Content
    [void]$rng.Find.Execute(" Структура", $true, $false, $false, $false, $false, $true, 1, $false, " Структура", 2)

    # Re-apply the red color to the single space between "приложения."
    # and "Схема" (this space used to carry that formatting before the
    # edit, and still should afterwards).
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute("приложения. Схема", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $spacePos = $rng2.Start + 11
        $spaceRng = $d.Range($spacePos, $spacePos + 1)
        $spaceRng.Font.Color = 188
    }
}
